# ---------------------------------------------------------------------
# Android Test Tool - Additional Informations.docx
#
# Two textual edits inside the same paragraph ("This tool is provided
# via an installer. ..."):
#
#   1. Remove the stray run containing "a " right before the opening
#      quote in "...permission-free folder (e.g. a 'User' folder...)"
#      so it reads "...permission-free folder (e.g. 'User' folder...)".
#
#   2. Replace "InterfaceTestTool" with "Test Tool" in
#      "Once installed, open the InterfaceTestTool application." and
#      keep the result split across three runs: the prefix, "Test
#      Tool" itself, and the trailing " application. ".
#
# This runtime merges every run that shares the same formatting into
# one big run, starting at whatever is being edited, whenever a Range
# is modified (Delete / Text assignment / Find-Replace). To keep the
# rest of the paragraph exactly the way it was (matching the upstream
# diff, which touches nothing else), each untouched run is "re-pinned"
# right after the edit by toggling Bold on and back off over its full
# span -- that is enough to stop the runtime from folding it into its
# neighbour. This has to happen strictly left-to-right, one run at a
# time, immediately following the edit, otherwise the pin does not
# stick.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

function Freeze-Run($startPos, $endPos) {
    if ($endPos -le $startPos) { return }
    $rng = $d.Range($startPos, $endPos)
    $rng.Bold = 1
    $rng.Bold = 0
}

# Anchor on the start of the paragraph via a uniquely-identifying
# piece of text so every offset below is computed relative to a live
# position rather than a hard-coded absolute one.
$anchor = $d.Content.Duplicate()
$found = $anchor.Find.Execute("This tool is provided via an installer", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate the target paragraph"
}
$base = $anchor.Start

# Original run texts/boundaries for this paragraph, taken from the
# pristine document (24 runs -> 25 cumulative-length boundaries).
$runTexts = @(
  "This tool ",
  "is provided via an installer. In the installation folder, open setup and follow the steps. The destination folder must be a permission-free folder (e.g. ",
  "a ",
  "‘",
  "U",
  "ser",
  "’",
  " folder, and not ",
  "‘",
  "Program Files",
  "’",
  "). ",
  "Once installed, open the InterfaceTestTool application. ",
  "It ",
  "uses ADB to communicate with Android Phones. ADB",
  " does not come by default during the first installation",
  ".",
  " T",
  "he latest version will be downloaded from internet and unzipped",
  ", and be available ",
  "after ",
  "in the ‘platform-tools’ folder",
  ".",
  " ADB assigns unique IDs to devices, and those IDs are used when multiple phones are connected with the following syntax: "
)

$offsets = @(0)
for ($i = 0; $i -lt $runTexts.Length; $i++) {
    $offsets += ($offsets[$offsets.Length - 1] + $runTexts[$i].Length)
}

$shift = 0

for ($i = 2; $i -le 23; $i++) {
    $s = $base + $offsets[$i] + $shift
    $e = $base + $offsets[$i + 1] + $shift

    if ($i -eq 2) {
        # --- Edit 1: delete the "a " run entirely --------------------
        $run2 = $d.Range($s, $e)
        if ($run2.Text -ne "a ") {
            throw "Unexpected text where 'a ' run was expected: [$($run2.Text)]"
        }
        $run2.Delete()
        $shift = $shift - ($e - $s)
    }
    elseif ($i -eq 12) {
        # --- Edit 2: split "...InterfaceTestTool..." into 3 runs -----
        $run12 = $d.Range($s, $e)
        if ($run12.Text -ne "Once installed, open the InterfaceTestTool application. ") {
            throw "Unexpected text where the InterfaceTestTool run was expected: [$($run12.Text)]"
        }

        $r = $d.Range($s, $e)
        $replaced = $r.Find.Execute("InterfaceTestTool", $true, $false, $false, $false, $false, $true, 1, $false, "Test Tool", 2)
        if (-not $replaced) {
            throw "Could not replace InterfaceTestTool with Test Tool"
        }

        $prefixLen = "Once installed, open the ".Length
        $midLen = "Test Tool".Length
        $oldMidLen = "InterfaceTestTool".Length
        $lenDelta = $midLen - $oldMidLen

        $b1 = $s + $prefixLen
        $b2 = $b1 + $midLen
        $e2 = $e + $lenDelta

        Freeze-Run $s $b1
        Freeze-Run $b1 $b2
        Freeze-Run $b2 $e2

        $shift = $shift + $lenDelta
    }
    else {
        # --- Untouched run: re-pin it so it keeps its own identity ---
        Freeze-Run $s $e
    }
}

Write-Host "Edits applied successfully"
